$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "59.566.93"
Set-TextValue $ws.Range("E2") "  +1.09%  "

Set-TextValue $ws.Range("D3") "2.602.64"
Set-TextValue $ws.Range("E3") "  +0.77%  "

Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.15%  "

Set-TextValue $ws.Range("D5") "539.58"
Set-TextValue $ws.Range("E5") "  +3.47%  "

Set-TextValue $ws.Range("D6") "141.77"
Set-TextValue $ws.Range("E6") "  +2.17%  "

Set-TextValue $ws.Range("E7") "  +0.09%  "

Set-TextValue $ws.Range("E8") "  +0.29%  "

Set-TextValue $ws.Range("E9") "  +0.25%  "

Set-TextValue $ws.Range("E10") "  +1.59%  "

Set-TextValue $ws.Range("D11") "0.335"
Set-TextValue $ws.Range("E11") "  +1.71%  "

Set-TextValue $ws.Range("E12") "  -0.66%  "

Set-TextValue $ws.Range("D13") "3.064.15"
Set-TextValue $ws.Range("E13") "  +0.63%  "

Set-TextValue $ws.Range("D14") "59.489.64"
Set-TextValue $ws.Range("E14") "  +1.03%  "

Set-TextValue $ws.Range("D15") "20.74"
Set-TextValue $ws.Range("E15") "  +1.08%  "

Set-TextValue $ws.Range("D16") "2.625.65"
Set-TextValue $ws.Range("E16") "  +1.39%  "

Set-TextValue $ws.Range("E17") "  +0.67%  "

Set-TextValue $ws.Range("D18") "341.56"
Set-TextValue $ws.Range("E18") "  +0.75%  "

Set-TextValue $ws.Range("E19") "  +1.82%  "

Set-TextValue $ws.Range("D20") "10.10"
Set-TextValue $ws.Range("E20") "  +0.49%  "

Set-TextValue $ws.Range("D21") "6.35"
Set-TextValue $ws.Range("E21") "  -1.34%  "

Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.06%  "

Set-TextValue $ws.Range("E23") "  +1.92%  "

Set-TextValue $ws.Range("E24") "  +1.72%  "

Set-TextValue $ws.Range("E25") "  -1.21%  "

Set-TextValue $ws.Range("D26") "0.993"
Set-TextValue $ws.Range("E26") "  -0.54%  "

Set-TextValue $ws.Range("D27") "7.24"
Set-TextValue $ws.Range("E27") "  +3.45%  "

Set-TextValue $ws.Range("D28") "0.0₃0745"
Set-TextValue $ws.Range("E28") "  +3.70%  "

Set-TextValue $ws.Range("E29") "  +0.02%  "

Set-TextValue $ws.Range("E30") "  +6.50%  "

Set-TextValue $ws.Range("D31") "5.85"
Set-TextValue $ws.Range("E31") "  -0.57%  "

Set-TextValue $ws.Range("D32") "18.82"
Set-TextValue $ws.Range("E32") "  +0.88%  "

Set-TextValue $ws.Range("D33") "149.58"
Set-TextValue $ws.Range("E33") "  +0.09%  "

Set-TextValue $ws.Range("D34") "4.01"
Set-TextValue $ws.Range("E34") "  +1.37%  "

Set-TextValue $ws.Range("E35") "  +0.77%  "

Set-TextValue $ws.Range("D36") "0.845"
Set-TextValue $ws.Range("E36") "  +4.40%  "

Set-TextValue $ws.Range("E37") "  -0.37%  "

Set-TextValue $ws.Range("E38") "  +0.62%  "

Set-TextValue $ws.Range("E39") "  +0.60%  "

Set-TextValue $ws.Range("E40") "  +0.21%  "

Set-TextValue $ws.Range("D41") "273.02"
Set-TextValue $ws.Range("E41") "  +0.62%  "

Set-TextValue $ws.Range("D42") "0.601"
Set-TextValue $ws.Range("E42") "  +1.55%  "

Set-TextValue $ws.Range("D43") "10.72"
Set-TextValue $ws.Range("E43") "  -0.12%  "

Set-TextValue $ws.Range("D44") "0.0950"
Set-TextValue $ws.Range("E44") "  +0.17%  "

Set-TextValue $ws.Range("E45") "  +1.49%  "

Set-TextValue $ws.Range("D46") "18.58"
Set-TextValue $ws.Range("E46") "  +4.24%  "

Set-TextValue $ws.Range("E47") "  +1.59%  "

Set-TextValue $ws.Range("D48") "1.941.61"
Set-TextValue $ws.Range("E48") "  -1.22%  "

Set-TextValue $ws.Range("E49") "  +0.28%  "

Set-TextValue $ws.Range("D50") "112.60"
Set-TextValue $ws.Range("E50") "  -1.09%  "

Set-TextValue $ws.Range("E51") "  +1.72%  "
